# refactor: remove vns and best descent
# Update the "C" column (gap-calculation baseline values) for rows 2-31.
# Columns D (per-row gap) and E (average gap) are formulas and will
# recalculate automatically from the new C values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    23067,
    23295,
    22887,
    22003,
    23072,
    23976,
    24326,
    22942,
    23192,
    23956,
    41124,
    40847,
    40503,
    43163,
    40707,
    41061,
    40528,
    42859,
    41539,
    41974,
    59412,
    62008,
    59299,
    59815,
    60043,
    58338,
    60718,
    60478,
    58802,
    59605
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}

$excel.Calculate()
